$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '41.748.32'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").Value = '2.217.20'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.15'

$ws.Range("E6").Value = '  -1.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.45'
$ws.Range("E7").Value = '  -2.52%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -4.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.37'
$ws.Range("E10").Value = '  -4.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  -1.92%  '

$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("E13").Value = '  -4.09%  '

$ws.Range("D14").Value = '2.548.22'
$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.16'
$ws.Range("E15").Value = '  -1.99%  '

$ws.Range("E16").Value = '  -2.81%  '

$ws.Range("D17").Value = '2.217.20'
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").Value = '41.662.37'
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("E19").Value = '  -4.92%  '

$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.75'
$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.90'
$ws.Range("E22").Value = '  +8.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '227.67'
$ws.Range("E23").Value = '  -1.53%  '

$ws.Range("E24").Value = '  -6.78%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.25'
$ws.Range("E26").Value = '  -4.77%  '

$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.86'
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.35'
$ws.Range("E31").Value = '  -3.52%  '

$ws.Range("E32").Value = '  -2.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.47'
$ws.Range("E33").Value = '  -5.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.45'
$ws.Range("E34").Value = '  +3.32%  '

$ws.Range("E35").Value = '  -1.46%  '

$ws.Range("E36").Value = '  -8.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  -4.12%  '

$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.95'

$ws.Range("E40").Value = '  -2.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '63.82'
$ws.Range("E41").Value = '  +1.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.59'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("E43").Value = '  -3.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.62'
$ws.Range("E44").Value = '  -2.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '101.59'
$ws.Range("E45").Value = '  -3.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.0993'
$ws.Range("E46").Value = '  -2.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  -1.17%  '

$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("E49").Value = '  -3.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.67'
$ws.Range("E50").Value = '  -1.59%  '

$ws.Range("D51").Value = '2.425.70'
$ws.Range("E51").Value = '  -0.55%  '
